$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: mark a LEVEL1 cell as "DONE" by pasting the format of an
# already-"DONE"-styled sibling cell (so fills/borders match exactly what
# Excel would reuse) and then writing the text.
function Set-DoneCell {
    param(
        [string]$TargetAddr,
        [string]$SourceAddr
    )
    $src = $ws.Range($SourceAddr)
    $dst = $ws.Range($TargetAddr)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $dst.Value = "DONE"
}

# L01r04 ("Rnd 04", LEVEL 1 column) terrain completion (excluding the
# fighting system block): mark the relevant LEVEL1 cells as DONE, matching
# the style already used by the analogous finished cells in the same rows.
Set-DoneCell "G6" "B6"
Set-DoneCell "B8" "G8"
Set-DoneCell "B9" "G9"
Set-DoneCell "B11" "B6"
Set-DoneCell "G11" "B6"
Set-DoneCell "G18" "B18"
Set-DoneCell "G19" "B19"
Set-DoneCell "B21" "B6"
Set-DoneCell "G21" "B6"

# "Complete terrain" / Rnd 04: LEVEL1 doesn't need the step, LEVEL2 is
# still in progress ("s").
$ws.Range("G16").Interior.Color = 0xA6A6A6
$ws.Range("G16").Value = "není třeba"
$ws.Range("H16").Value = "s"

# View table: move the active selection.
$ws.Range("H12").Select() | Out-Null
